$wb = $excel.ActiveWorkbook

# --- Sheet "Đơn phụ phẫu 1": insert a new data row (708 / Nguyễn Ngọc Hân) ---
# before row 17, pushing the existing "Tổng" row down to row 18.
$ws1 = $wb.Worksheets.Item("Đơn phụ phẫu 1")
$ws1.Rows.Item(17).Insert()

$ws1.Cells.Item(17, 1).Value = "HD-LUXURY"
$ws1.Cells.Item(17, 2).Value = 708

# Column C holds dates as plain text (e.g. "08-01-2024") rather than real
# date values elsewhere in the sheet, so force text formatting, write the
# text, then drop the number-format override again so the cell ends up
# with the sheet's default (General) style, matching its neighbours.
$ws1.Cells.Item(17, 3).NumberFormat = "@"
$ws1.Cells.Item(17, 3).Value = "08-31-2024"
$ws1.Cells.Item(17, 3).ClearFormats()

$ws1.Cells.Item(17, 4).Value = "CẦN THƠ"
$ws1.Cells.Item(17, 5).Value = "Nguyễn Ngọc Hân"
$ws1.Cells.Item(17, 6).Value = "Cá nhân"
$ws1.Cells.Item(17, 7).Value = "Nâng mũi"
$ws1.Cells.Item(17, 8).Value = "Lâm Hoàng Phú"
$ws1.Cells.Item(17, 9).Value = 100000

# Update the "Tổng" row (now row 18) totals.
$ws1.Cells.Item(18, 2).Value = 16
$ws1.Cells.Item(18, 9).Value = 1200000

# --- Sheet "Lương": refresh the computed totals ---
$ws5 = $wb.Worksheets.Item("Lương")
$ws5.Cells.Item(1, 2).Value = 13
$ws5.Cells.Item(2, 2).Value = 27
$ws5.Cells.Item(3, 2).Value = 3857142.857142857
$ws5.Cells.Item(8, 2).Value = 1200000
$ws5.Cells.Item(35, 2).Value = 1677142.857142857
$ws5.Cells.Item(38, 2).Value = 1777142.857142857
